# Combine two similar phrases into one on the "list" sheet.
#
# Row 10 (A10): "Echo / feedback / cutting out"  -> "echo / feedback"
# Row 54 (A54): "You cut out"                     -> "You're cutting out"
#
# After the edit, select A54 (matching the new author's last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

$ws.Range("A10").Value = "echo / feedback"
$ws.Range("A54").Value = "You're cutting out"

$ws.Activate()
$ws.Range("A54").Select()
